$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data read-in used to skip a row (leaving a bogus all-zero "middle" row in
# the sheet); now it does not, so that spurious row must be removed and every
# row below it shifts up by one.
$ws.Rows.Item(32).Delete() | Out-Null

# Reflect where the cursor ended up after the edit.
$ws.Range("Q41").Select() | Out-Null
